# Insert a new row at position 143 (pushes existing rows 143-186 down to 144-187)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(143).Insert()

# Populate the newly inserted row 143 with the new record's data.
$ws.Cells.Item(143, 1).Value = 5
$ws.Cells.Item(143, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(143, 3).Value = "Maule"
$ws.Cells.Item(143, 4).Value = 44463
$ws.Cells.Item(143, 5).Value = 7
$ws.Cells.Item(143, 6).Value = 100114013
$ws.Cells.Item(143, 7).Value = "Zanahoria"
$ws.Cells.Item(143, 8).Value = "Sin especificar"
$ws.Cells.Item(143, 9).Value = "Primera"
$ws.Cells.Item(143, 10).Value = 400
$ws.Cells.Item(143, 11).Value = 6000
$ws.Cells.Item(143, 12).Value = 6000
$ws.Cells.Item(143, 13).Value = 6000
$ws.Cells.Item(143, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(143, 15).Value = "Región de Ñuble"
$ws.Cells.Item(143, 16).Value = 300
$ws.Cells.Item(143, 17).Value = 20
$ws.Cells.Item(143, 18).Value = "Hortaliza"
